$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for Coliflor (Vega Modelo de Temuco).
# It belongs right above the current row 223, so push everything at/after
# row 223 down by one (this also preserves the date-formatted style on
# column D for every shifted row) and then fill in the new record's values.
$ws.Rows.Item(223).Insert()

$ws.Range("A223").Value = 10
$ws.Range("B223").Value = "Vega Modelo de Temuco"
$ws.Range("C223").Value = "La Araucanía"
$ws.Range("D223").Value = 44505
$ws.Range("E223").Value = 9
$ws.Range("F223").Value = 100112008
$ws.Range("G223").Value = "Coliflor"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 2150
$ws.Range("K223").Value = 800
$ws.Range("L223").Value = 900
$ws.Range("M223").Value = 842
$ws.Range("N223").Value = "$/unidad"
$ws.Range("O223").Value = "Región del Maule"
$ws.Range("P223").Value = 842
$ws.Range("Q223").Value = 1
$ws.Range("R223").Value = "Hortaliza"
